$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")

# --- Fill in Date (E) and Version (F) for rows 341-366 ---
# E column: date serial 44379 (2021-07-02), already formatted via existing style
# F column: version text "1.0000.0000000" (stored as text)
for ($r = 341; $r -le 366; $r++) {
    $ws.Cells.Item($r, 5).Value = 44379
    $ws.Cells.Item($r, 6).Value = "1.0000.0000000"
}

# --- Update Description (C) for rows 353-366 with the new, more specific API names ---
$ws.Cells.Item(353, 3).Value = "Membatalkan Penghapusan Data Jenis Barang"
$ws.Cells.Item(354, 3).Value = "Membatalkan Penghapusan Data Institusi"
$ws.Cells.Item(355, 3).Value = "Membatalkan Penghapusan Data Cabang Institusi"
$ws.Cells.Item(356, 3).Value = "Membatalkan Penghapusan Data Periode"
$ws.Cells.Item(357, 3).Value = "Membatalkan Penghapusan Data Orang"
$ws.Cells.Item(358, 3).Value = "Membatalkan Penghapusan Data Akun E-Mail Orang"
$ws.Cells.Item(359, 3).Value = "Membatalkan Penghapusan Data Akun Sosial Media Orang"
$ws.Cells.Item(360, 3).Value = "Membatalkan Penghapusan Data Jenis Kelamin Orang"
$ws.Cells.Item(361, 3).Value = "Membatalkan Penghapusan Data Produk"
$ws.Cells.Item(362, 3).Value = "Membatalkan Penghapusan Data Jenis Produk"
$ws.Cells.Item(363, 3).Value = "Membatalkan Penghapusan Data Unit Kuantitas"
$ws.Cells.Item(364, 3).Value = "Membatalkan Penghapusan Data Agama"
$ws.Cells.Item(365, 3).Value = "Membatalkan Penghapusan Data Media Sosial"
$ws.Cells.Item(366, 3).Value = "Membatalkan Penghapusan Data Merk Dagang"

# --- Update frozen-pane scroll position / active selection to reflect the new bottom of the list ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 359
$ws.Range("C364").Select()
